$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells keep their original Text storage (inline string) rather than
# being auto-converted to numbers/percentages by Excel when we assign numeric-
# or percent-looking strings. Pre-formatting as Text ("@") before the assignment
# mirrors how these price/volume columns were authored as text in the workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.19%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.18%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.134"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.05%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07636"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.690"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.53%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9350"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.15%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1256"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.52%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.63%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09042"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04166"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.20%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1056"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.62%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001267"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.83%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005825"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.59%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.364"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.41%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.47%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3360"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.73%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.406"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "21.14%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.10%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04038"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.30%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001264"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.19%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004052"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.46%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.47%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02483"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.34%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05197"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.51%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1300"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007366"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.43%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002168"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "17.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008163"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.31%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3136"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.15%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006656"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.15%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.46%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2707"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "59.37%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004215"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "2.91%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002108"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.46%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.46%"
